# Update "Förändrad" (column C) dates for rows 2-11 from 2023-10-05 (45204)
# to 2023-10-08 (45207), as reflected in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value = 45207
    }
}
